# Update "want to go" (想去人数) counts in column F for the rows that
# changed between scrapes. The same underlying data is duplicated on the
# "展览" sheet and the "全部类型" sheet, so both need the same update.

$wb = $excel.ActiveWorkbook

$updates = @{
    3 = 1864
    6 = 1122
    7 = 50
    8 = 5958
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
